$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2633.6843
$ws.Range("I19").Value = 1960.5333
$ws.Range("J19").Value = 5158
$ws.Range("K19").Value = 1960.5333
$ws.Range("L19").Value = 5158
$ws.Range("M19").Value = -1785.5333
$ws.Range("N19").Value = -5508
$ws.Range("H33").Value = 304.2
$ws.Range("I33").Value = 252.05882
$ws.Range("K33").Value = 252.05882
$ws.Range("M33").Value = -23.05882
$ws.Range("H88").Value = 15908194
$ws.Range("J88").Value = 51134.855
$ws.Range("L88").Value = 51134.855
$ws.Range("N88").Value = -51946.855
$ws.Range("H91").Value = 15908194
$ws.Range("J91").Value = 51134.855
$ws.Range("L91").Value = 51134.855
$ws.Range("N91").Value = -53942.855
$ws.Range("H132").Value = 1925.1538
$ws.Range("I132").Value = 1593.9584
$ws.Range("K132").Value = 4781.8752
$ws.Range("M132").Value = -2251.8752
$ws.Range("H137").Value = 2159.4614
$ws.Range("I137").Value = 1632.64
$ws.Range("K137").Value = 4897.92
$ws.Range("M137").Value = -2347.92
$ws.Range("H138").Value = 4049.1353
$ws.Range("I138").Value = 1055.1212
$ws.Range("K138").Value = 3165.3636
$ws.Range("M138").Value = 1974.6364
$ws.Range("H141").Value = 2200.5715
$ws.Range("I141").Value = 2097.037
$ws.Range("K141").Value = 6291.110999999999
$ws.Range("M141").Value = -1111.110999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2608.7407
$ws.Range("I2").Value = 1330.1177
$ws.Range("J2").Value = 4782.4
$ws.Range("K2").Value = 1330.1177
$ws.Range("L2").Value = 4782.4
$ws.Range("M2").Value = -1217.1177
$ws.Range("N2").Value = -5008.4
$ws.Range("H32").Value = 1627612.5
$ws.Range("I32").Value = 1764137.5
$ws.Range("K32").Value = 1764137.5
$ws.Range("M32").Value = -1763850.5
$ws.Range("H45").Value = 6046.4546
$ws.Range("I45").Value = 1757.909
$ws.Range("K45").Value = 1757.909
$ws.Range("M45").Value = -1380.909
$ws.Range("H61").Value = 7898.393
$ws.Range("I61").Value = 3824.6155
$ws.Range("K61").Value = 3824.6155
$ws.Range("M61").Value = -3612.6155
$ws.Range("H74").Value = 2160.3965
$ws.Range("J74").Value = 3930.0454
$ws.Range("L74").Value = 3930.0454
$ws.Range("N74").Value = -5678.0454
$ws.Range("H77").Value = 2160.3965
$ws.Range("J77").Value = 3930.0454
$ws.Range("L77").Value = 19650.227
$ws.Range("N77").Value = -28386.227
$ws.Range("H102").Value = 1557.4286
$ws.Range("I102").Value = 1500.3334
$ws.Range("J102").Value = 1900
$ws.Range("K102").Value = 1500.3334
$ws.Range("L102").Value = 1900
$ws.Range("M102").Value = 121.6666
$ws.Range("N102").Value = -5144
$ws.Range("H116").Value = 2608.7407
$ws.Range("I116").Value = 1330.1177
$ws.Range("J116").Value = 4782.4
$ws.Range("K116").Value = 1330.1177
$ws.Range("L116").Value = 4782.4
$ws.Range("M116").Value = 963.8823
$ws.Range("N116").Value = -9370.4
$ws.Range("H132").Value = 3408.3013
$ws.Range("I132").Value = 2028.2549
$ws.Range("K132").Value = 6084.7647
$ws.Range("M132").Value = -3554.7647
$ws.Range("H135").Value = 70213.5
$ws.Range("J135").Value = 70213.5
$ws.Range("L135").Value = 70213.5
$ws.Range("N135").Value = -80353.5
$ws.Range("H136").Value = 7898.393
$ws.Range("I136").Value = 3824.6155
$ws.Range("K136").Value = 11473.8465
$ws.Range("M136").Value = -8923.8465

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2608.7407
$ws.Range("I3").Value = 1330.1177
$ws.Range("J3").Value = 4782.4
$ws.Range("K3").Value = 1330.1177
$ws.Range("L3").Value = 4782.4
$ws.Range("M3").Value = -1216.1177
$ws.Range("N3").Value = -5010.4
$ws.Range("H22").Value = 250.77777
$ws.Range("I22").Value = 250.77777
$ws.Range("K22").Value = 250.77777
$ws.Range("M22").Value = -77.77777
$ws.Range("H105").Value = 2611.342
$ws.Range("I105").Value = 2335.8965
$ws.Range("J105").Value = 3498.889
$ws.Range("K105").Value = 2335.8965
$ws.Range("L105").Value = 3498.889
$ws.Range("M105").Value = -588.8964999999998
$ws.Range("N105").Value = -6992.889

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 5010
$ws.Range("J23").Value = 5010
$ws.Range("L23").Value = 5010
$ws.Range("N23").Value = -5490
$ws.Range("H27").Value = 5010
$ws.Range("J27").Value = 5010
$ws.Range("L27").Value = 5010
$ws.Range("N27").Value = -5394
$ws.Range("H93").Value = 53942
$ws.Range("J93").Value = 53942
$ws.Range("L93").Value = 53942
$ws.Range("N93").Value = -57686
$ws.Range("H99").Value = 4012.48
$ws.Range("I99").Value = 2705.4443
$ws.Range("J99").Value = 7373.4287
$ws.Range("K99").Value = 2705.4443
$ws.Range("L99").Value = 7373.4287
$ws.Range("M99").Value = -1207.4443
$ws.Range("N99").Value = -10369.4287
$ws.Range("H107").Value = 1338.0278
$ws.Range("I107").Value = 1005.26086
$ws.Range("J107").Value = 1926.7693
$ws.Range("K107").Value = 1005.26086
$ws.Range("L107").Value = 1926.7693
$ws.Range("M107").Value = 914.73914
$ws.Range("N107").Value = -5766.7693
$ws.Range("H126").Value = 4012.48
$ws.Range("I126").Value = 2705.4443
$ws.Range("J126").Value = 7373.4287
$ws.Range("K126").Value = 8116.3329
$ws.Range("L126").Value = 22120.2861
$ws.Range("M126").Value = -5646.3329
$ws.Range("N126").Value = -27060.2861
$ws.Range("H132").Value = 5680.4736
$ws.Range("I132").Value = 2309.7334
$ws.Range("K132").Value = 6929.2002
$ws.Range("M132").Value = -4399.2002
$ws.Range("H134").Value = 4639.5576
$ws.Range("I134").Value = 1358.2858
$ws.Range("K134").Value = 4074.8574
$ws.Range("M134").Value = -1539.8574

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 606.2
$ws.Range("I114").Value = 223.14285
$ws.Range("J114").Value = 1500
$ws.Range("K114").Value = 669.4285500000001
$ws.Range("L114").Value = 4500
$ws.Range("M114").Value = 2584.57145
$ws.Range("N114").Value = -11008

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 616693.1
$ws.Range("I107").Value = 1335086.5
$ws.Range("K107").Value = 1335086.5
$ws.Range("M107").Value = -1333166.5
$ws.Range("H122").Value = 3574292.5
$ws.Range("I122").Value = 4764141
$ws.Range("K122").Value = 14292423
$ws.Range("M122").Value = -14289973
$ws.Range("H132").Value = 4624.125
$ws.Range("I132").Value = 1853.0588
$ws.Range("K132").Value = 5559.1764
$ws.Range("M132").Value = -3029.1764

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3868.0967
$ws.Range("I7").Value = 2086.682
$ws.Range("K7").Value = 2086.682
$ws.Range("M7").Value = -1974.682
$ws.Range("H61").Value = 3587.9697
$ws.Range("I61").Value = 1619
$ws.Range("K61").Value = 1619
$ws.Range("M61").Value = -1417
$ws.Range("H100").Value = 3790.2942
$ws.Range("J100").Value = 5750.8335
$ws.Range("L100").Value = 5750.8335
$ws.Range("N100").Value = -6832.8335
$ws.Range("H113").Value = 3587.9697
$ws.Range("I113").Value = 1619
$ws.Range("K113").Value = 1619
$ws.Range("M113").Value = 551
$ws.Range("H126").Value = 3868.0967
$ws.Range("I126").Value = 2086.682
$ws.Range("K126").Value = 6260.045999999999
$ws.Range("M126").Value = -3790.045999999999
$ws.Range("H132").Value = 10644839
$ws.Range("I132").Value = 21741498
$ws.Range("K132").Value = 65224494
$ws.Range("M132").Value = -65221964

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 28582386
$ws.Range("I81").Value = 3424.75
$ws.Range("J81").Value = 66687668
$ws.Range("K81").Value = 6849.5
$ws.Range("L81").Value = 133375336
$ws.Range("M81").Value = -5788.5
$ws.Range("N81").Value = -133377458
$ws.Range("H84").Value = 28582386
$ws.Range("I84").Value = 3424.75
$ws.Range("J84").Value = 66687668
$ws.Range("K84").Value = 34247.5
$ws.Range("L84").Value = 666876680
$ws.Range("M84").Value = -28943.5
$ws.Range("N84").Value = -666887288
$ws.Range("H136").Value = 23838044
$ws.Range("I136").Value = 50000820
$ws.Range("K136").Value = 150002460
$ws.Range("M136").Value = -149999910
